# "drag drop input ok" - update time-estimate inputs (column B) that were
# re-entered via drag & drop on the sheet; dependent formulas in column C
# recalc automatically. Also moves the active selection to where the user
# last clicked (B13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B5").Value = 1
$ws.Range("B6").Value = 1.5
$ws.Range("B12").Value = 1.5

$ws.Activate()
$ws.Range("B13").Select()
